$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the crypto price/volume table (columns D = Price, E = Volume(1h)).
# D-column values that look like plain numbers are written with a leading
# apostrophe so Excel keeps them as literal text (matching the source
# formatting, e.g. trailing zeros like "588.40" or "1.00") instead of
# coercing them into numeric cells.
$ws.Range("D2").Value = "67.911.45"
$ws.Range("E2").Value = "  +1.10%  "
$ws.Range("D3").Value = "2.487.20"
$ws.Range("E3").Value = "  +0.19%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'588.40"
$ws.Range("E5").Value = "  +0.54%  "
$ws.Range("D6").Value = "'174.50"
$ws.Range("E6").Value = "  +1.04%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "'0.514"
$ws.Range("E8").Value = "  -0.26%  "
$ws.Range("D9").Value = "'0.144"
$ws.Range("E9").Value = "  +4.13%  "
$ws.Range("E10").Value = "  -1.73%  "
$ws.Range("D11").Value = "'4.97"
$ws.Range("E11").Value = "  +0.49%  "
$ws.Range("E12").Value = "  +0.31%  "
$ws.Range("D13").Value = "2.940.45"
$ws.Range("E13").Value = "  +0.29%  "
$ws.Range("D14").Value = "'25.33"
$ws.Range("E14").Value = "  -0.84%  "
$ws.Range("D15").Value = "67.828.98"
$ws.Range("E15").Value = "  +1.14%  "
$ws.Range("D16").Value = "'0.0000170"
$ws.Range("E16").Value = "  -0.43%  "
$ws.Range("D17").Value = "2.483.76"
$ws.Range("E17").Value = "  -1.93%  "
$ws.Range("D18").Value = "'10.83"
$ws.Range("E18").Value = "  -1.29%  "
$ws.Range("D19").Value = "'7.40"
$ws.Range("E19").Value = "  -2.22%  "
$ws.Range("D20").Value = "'346.88"
$ws.Range("E20").Value = "  -1.00%  "
$ws.Range("D21").Value = "'4.12"
$ws.Range("E21").Value = "  +2.13%  "
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("D23").Value = "'70.78"
$ws.Range("E23").Value = "  +2.59%  "
$ws.Range("E24").Value = "  -1.99%  "
$ws.Range("D25").Value = "'1.69"
$ws.Range("E25").Value = "  -7.18%  "
$ws.Range("D26").Value = "'8.84"
$ws.Range("E26").Value = "  -3.69%  "
$ws.Range("E27").Value = "  +0.10%  "
$ws.Range("E28").Value = "  -0.09%  "
$ws.Range("D29").Value = "0.0₃0889"
$ws.Range("D30").Value = "'498.88"
$ws.Range("E30").Value = "  -1.49%  "
$ws.Range("D31").Value = "'7.73"
$ws.Range("E31").Value = "  +0.22%  "
$ws.Range("D32").Value = "'1.25"
$ws.Range("E32").Value = "  +0.01%  "
$ws.Range("E33").Value = "  -0.46%  "
$ws.Range("D34").Value = "'1.00"
$ws.Range("E34").Value = "  +0.04%  "
$ws.Range("D35").Value = "'164.60"
$ws.Range("E35").Value = "  +1.43%  "
$ws.Range("E36").Value = "  +1.47%  "
$ws.Range("D38").Value = "'18.28"
$ws.Range("E38").Value = "  +0.71%  "
$ws.Range("E39").Value = "  -0.03%  "
$ws.Range("D40").Value = "'1.31"
$ws.Range("E40").Value = "  -2.36%  "
$ws.Range("E41").Value = "  +2.08%  "
$ws.Range("E42").Value = "  -1.61%  "
$ws.Range("D43").Value = "'4.77"
$ws.Range("E43").Value = "  -1.28%  "
$ws.Range("D44").Value = "'2.38"
$ws.Range("E44").Value = "  -0.11%  "
$ws.Range("D45").Value = "'147.71"
$ws.Range("E45").Value = "  +2.89%  "
$ws.Range("D46").Value = "'3.52"
$ws.Range("E46").Value = "  +1.16%  "
$ws.Range("D47").Value = "'0.511"
$ws.Range("E47").Value = "  -0.96%  "
$ws.Range("D48").Value = "0.0₆0253"
$ws.Range("E48").Value = "  -3.79%  "
$ws.Range("E49").Value = "  -0.37%  "
$ws.Range("D50").Value = "'1.56"
$ws.Range("E50").Value = "  -1.57%  "
$ws.Range("D51").Value = "'0.576"
$ws.Range("E51").Value = "  -1.34%  "
